# Publish terminology IG 2.0.0
# - Metadata sheet: bump Version/Date/Count
# - Rename old "Concepts" sheet to "Properties" and replace its content with
#   the CodeSystem property definitions (status / effectiveDate)
# - Add a brand-new "Concepts" sheet with the full (now 27-row) concept list

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Metadata sheet updates
# ---------------------------------------------------------------------------
$meta = $wb.Worksheets.Item(1)
$meta.Range("B3").Value = "2.0.0"

# "2025-09-22" looks like a date to the engine's smart-typing, so force it to
# stay plain text (like the original "2022-09-01" string) and then restore
# the untouched general-text formatting that the rest of column B uses.
$meta.Range("B8").NumberFormat = "@"
$meta.Range("B8").Value = "2025-09-22"
$meta.Range("A8").Copy()
$meta.Range("B8").PasteSpecial(-4122)

$meta.Range("B22").Value = "27"

# ---------------------------------------------------------------------------
# 2. Rename the existing "Concepts" sheet to "Properties" and rewrite it
# ---------------------------------------------------------------------------
$props = $wb.Worksheets.Item(2)
$props.Name = "Properties"

# Drop the rows that belonged to the old concept table (21 rows -> 3 rows)
$props.Range("A4:D21").Clear()

$propsHeader = New-Object 'object[,]' 1,4
$propsHeader[0,0] = "Code"
$propsHeader[0,1] = "Uri"
$propsHeader[0,2] = "Description"
$propsHeader[0,3] = "Type"
$props.Range("A1:D1").Value = $propsHeader

$propsBody = New-Object 'object[,]' 2,4
$propsBody[0,0] = "status"
$propsBody[0,1] = "http://hl7.org/fhir/concept-properties#status"
$propsBody[0,2] = "A property that indicates the status of the concept. One of active, experimental, deprecated, or retired."
$propsBody[0,3] = "code"
$propsBody[1,0] = "effectiveDate"
$propsBody[1,1] = "http://hl7.org/fhir/concept-properties#effectiveDate"
$propsBody[1,2] = "The date at which the concept status was last changed."
$propsBody[1,3] = "dateTime"
$props.Range("A2:D3").Value = $propsBody

# ---------------------------------------------------------------------------
# 3. Insert a fresh "Concepts" sheet (right after "Properties") with the
#    updated 27-row concept list
# ---------------------------------------------------------------------------
$concepts = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $props)
$concepts.Name = "Concepts"

$data = New-Object 'object[,]' 28,4
$data[0,0] = "Level"
$data[0,1] = "Code"
$data[0,2] = "Display"
$data[0,3] = "Definition"
$data[1,0] = "1"
$data[1,1] = "MCM:FDIS20#AcuteAmbulant"
$data[1,2] = "Start hospital stay-acute ambulant"
$data[1,3] = ""
$data[2,0] = "1"
$data[2,1] = "MCM:FDIS20#AdmissionInpatient"
$data[2,2] = "Start hospital stay-admission"
$data[2,3] = ""
$data[3,0] = "1"
$data[3,1] = "MCM:FDIS20#OnLeave"
$data[3,2] = "Start leave"
$data[3,3] = ""
$data[4,0] = "1"
$data[4,1] = "MCM:FDIS20#EndOfLeave"
$data[4,2] = "end leave"
$data[4,3] = ""
$data[5,0] = "1"
$data[5,1] = "MCM:FDIS20#EndHospitalStay"
$data[5,2] = "End hospital Stay - patient discharged to home/primary sector"
$data[5,3] = ""
$data[6,0] = "1"
$data[6,1] = "MCM:FDIS91#outpatient"
$data[6,2] = "Outpatient"
$data[6,3] = ""
$data[7,0] = "1"
$data[7,1] = "MCM:FDIS91#decease"
$data[7,2] = "Decease"
$data[7,3] = ""
$data[8,0] = "1"
$data[8,1] = "MCM:FDIS91#carecoordination"
$data[8,2] = "Care Coordination"
$data[8,3] = ""
$data[9,0] = "1"
$data[9,1] = "MCM:FDIS91#assistive-devices"
$data[9,2] = "Assistive Devices"
$data[9,3] = ""
$data[10,0] = "1"
$data[10,1] = "MCM:FDIS91#medicine"
$data[10,2] = "Medicine"
$data[10,3] = ""
$data[11,0] = "1"
$data[11,1] = "MCM:FDIS91#psychiatry-social-disability"
$data[11,2] = "Psychiatry, Social, Disability"
$data[11,3] = ""
$data[12,0] = "1"
$data[12,1] = "MCM:FDIS91#healthcare"
$data[12,2] = "Healthcare"
$data[12,3] = ""
$data[13,0] = "1"
$data[13,1] = "MCM:FDIS91#nursing"
$data[13,2] = "Nursing"
$data[13,3] = ""
$data[14,0] = "1"
$data[14,1] = "MCM:FDIS91#telemedicine"
$data[14,2] = "Telemedicine"
$data[14,3] = ""
$data[15,0] = "1"
$data[15,1] = "MCM:FDIS91#training"
$data[15,2] = "Training"
$data[15,3] = ""
$data[16,0] = "1"
$data[16,1] = "MCM:FDIS91#discharge"
$data[16,2] = "Discharge"
$data[16,3] = ""
$data[17,0] = "1"
$data[17,1] = "MCM:FDIS91#regarding-referral"
$data[17,2] = "Regarding Referral"
$data[17,3] = ""
$data[18,0] = "1"
$data[18,1] = "MCM:FDIS91#assessment"
$data[18,2] = "Assessment"
$data[18,3] = ""
$data[19,0] = "1"
$data[19,1] = "MCM:FDIS91#home-care-assessment"
$data[19,2] = "Home Care Assessment"
$data[19,3] = ""
$data[20,0] = "1"
$data[20,1] = "MCM:FDIS91#examination-results"
$data[20,2] = "Examination Results"
$data[20,3] = ""
$data[21,0] = "1"
$data[21,1] = "MCM:FDIS91#alcohol-and-drug-treatment"
$data[21,2] = "Alcohol and Drug Treatment"
$data[21,3] = ""
$data[22,0] = "1"
$data[22,1] = "MCM:FDIS91#acute-ambulant"
$data[22,2] = "Acute Ambulant"
$data[22,3] = ""
$data[23,0] = "1"
$data[23,1] = "MCM:FDIS91#extended-care-responsibility"
$data[23,2] = "Extended Care Responsibility"
$data[23,3] = ""
$data[24,0] = "1"
$data[24,1] = "MCM:FDIS91#other"
$data[24,2] = "Other"
$data[24,3] = ""
$data[25,0] = "1"
$data[25,1] = "MCM:FCTL#ok"
$data[25,2] = "OK"
$data[25,3] = ""
$data[26,0] = "1"
$data[26,1] = "MCM:FCTL#transient-error"
$data[26,2] = "Transient Error"
$data[26,3] = ""
$data[27,0] = "1"
$data[27,1] = "MCM:FCTL#fatal-error"
$data[27,2] = "Fatal Error"
$data[27,3] = ""
$concepts.Range("A1:D28").Value = $data

# Copy the header / body formatting from the Properties sheet onto the new
# Concepts sheet so the look (bold header row, wrapped body rows) matches.
$props.Range("A1").Copy()
$concepts.Range("A1:D1").PasteSpecial(-4122)
$props.Range("A2").Copy()
$concepts.Range("A2:D28").PasteSpecial(-4122)

$meta.Select()
